# Append the newest risk run (80d5083a-00d0-4303-8126-8a0b0cc7ecc3) to both
# the "Risks" summary sheet and the "HazardDetails" sheet, as produced by
# saving the whylogs visualization run's release artifacts.

$wb = $excel.ActiveWorkbook

# --- Risks sheet: append row 6 ---
$risks = $wb.Worksheets.Item("Risks")

$risks.Cells.Item(6, 1).Value = "80d5083a-00d0-4303-8126-8a0b0cc7ecc3"
$risks.Cells.Item(6, 2).Value = "2025-05-18T21:50:58.951381"
$risks.Cells.Item(6, 3).Value = 0.525
$risks.Cells.Item(6, 4).Value = 0.25
$risks.Cells.Item(6, 5).Value = 0.8
$risks.Cells.Item(6, 6).Value = "HIGH"
$risks.Cells.Item(6, 7).Value = "PENDING"
$risks.Cells.Item(6, 8).Value = "Unfair bias against protected demographic groups"
$risks.Cells.Item(6, 9).Value = "Re-sample training data; add fairness constraints or post-processing techniques"

# --- HazardDetails sheet: append row 4 ---
$hazards = $wb.Worksheets.Item("HazardDetails")

$hazards.Cells.Item(4, 1).Value = "80d5083a-00d0-4303-8126-8a0b0cc7ecc3"
$hazards.Cells.Item(4, 2).Value = "2025-05-18T21:50:58.951381"
$hazards.Cells.Item(4, 3).Value = 0.525
$hazards.Cells.Item(4, 4).Value = "bias_protected_groups"
$hazards.Cells.Item(4, 5).Value = "Unfair bias against protected demographic groups"
$hazards.Cells.Item(4, 6).Value = "HIGH"
$hazards.Cells.Item(4, 7).Value = "Re-sample training data; add fairness constraints or post-processing techniques"
$hazards.Cells.Item(4, 8).Value = "num__AGE_YEARS: 1.000 disparity`n"
